$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters corresponding to each value column in $data rows (B, D, E, F, G, L, M, N, O)
$cols = @(2, 4, 5, 6, 7, 12, 13, 14, 15)

$data = @(
    @(16.39862850593824, 3.477402772200163, 22.80869768991504, 22.0281974001364, 3.590684002881349, 8.577136736468916, 15.07501039886376, 19.44456028645104, 19.21709317571765),
    @(16.23145348657, 3.467115742155728, 22.65606221959466, 21.77035295446004, 3.593229942687646, 8.477002520740809, 14.98869745960153, 19.46701003805822, 19.09626238823792),
    @(16.13081926170712, 3.460731157103996, 22.56257219661076, 21.61818102675239, 3.594877293711931, 8.415957553353433, 14.93751499042491, 19.48322876891323, 19.02783396003856),
    @(16.09035770080796, 3.458112847975432, 22.52455540483771, 21.55779249655497, 3.595569829073876, 8.391215761860138, 14.91712957496795, 19.49045166775327, 19.00142301135808),
    @(16.08367337119525, 3.457677099386604, 22.51824833265212, 21.54786522316435, 3.59568610813754, 8.387116241408096, 14.9137735375829, 19.49168812818974, 18.99712721414814),
    @(16.13027131085902, 3.460695911807268, 22.56205912921217, 21.61735993607413, 3.594886547447998, 8.415623300084391, 14.93723813480069, 19.48332369298637, 19.02747177256051),
    @(16.34059434050399, 3.473870363183526, 22.75603572680652, 21.9380647098966, 3.591544423375593, 8.542534286815487, 15.04488280442568, 19.45179657413881, 19.17425077079178),
    @(16.76712503867862, 3.499135974694868, 23.13719644820755, 22.61215259341746, 3.585654888474467, 8.793802973789456, 15.26965065511827, 19.40922332261857, 19.50656426909281),
    @(17.08653557026863, 3.517315441604726, 23.41630852782912, 23.12997416579037, 3.581728385601701, 8.978431999761465, 15.44213140974174, 19.38959150105454, 19.77599759646522),
    @(17.23260890836242, 3.525494869523125, 23.54278474989312, 23.36935490805285, 3.580028139315175, 9.062128193396807, 15.52197109801086, 19.38316917464588, 19.90362350031512),
    @(17.28798855837017, 3.52857851608574, 23.59058196243833, 23.46045988280879, 3.579396585773896, 9.093755037552191, 15.55238319605516, 19.38109605307879, 19.95264031095401),
    @(17.27605940464187, 3.527915019391897, 23.58029267987206, 23.44081995585492, 3.57953205643497, 9.086947006802943, 15.54582577923206, 19.381526605581, 19.9420537304818),
    @(17.23716398824407, 3.525748842134089, 23.54671913741425, 23.37684154581783, 3.5799759350303, 9.06473163872592, 15.52446965404133, 19.38299143597933, 19.90764259482213),
    @(17.21334650721995, 3.524420187455372, 23.52614103623102, 23.33770965223301, 3.580249422449693, 9.051114614574187, 15.51141107860569, 19.38393536774366, 19.88665312472102),
    @(17.07700090982118, 3.516779014850398, 23.4080308864786, 23.1143990186506, 3.581841224436557, 8.972954218685517, 15.43693975113373, 19.39006153495293, 19.76775537577109),
    @(16.99352201230311, 3.512067717541475, 23.33542949291475, 22.97831770457564, 3.5828397080552, 8.924912451897088, 15.39159296785732, 19.39446090557552, 19.69608338933148),
    @(16.9455825485659, 3.509349427993227, 23.29362618999944, 22.90041253799333, 3.583422102327045, 8.897254033992207, 15.36564152458357, 19.39722751237408, 19.65533825786925),
    @(16.92936543247507, 3.508427627761733, 23.27946538243044, 22.87410074645416, 3.583620682979216, 8.887885664230105, 15.35687787192555, 19.39820486042964, 19.64162606373893),
    @(17.00240101677, 3.512570125549615, 23.34316286877331, 22.99276666908855, 3.582732580656971, 8.930029460660634, 15.39640681080007, 19.39396815161384, 19.70366373512421),
    @(17.24858711999805, 3.526385479075548, 23.55658332561291, 23.39562189784542, 3.579845224021257, 9.071258848152874, 15.5307377768732, 19.38255145488117, 19.91773164180891),
    @(17.40984097692695, 3.535334279309572, 23.6954911841687, 23.66153051485522, 3.578029791639091, 9.163158201426699, 15.61956229475549, 19.37718088080158, 20.0616260920468),
    @(17.32375918595839, 3.530565721786357, 23.62141446141613, 23.5194008651951, 3.578992190053821, 9.11415470275225, 15.57206710487862, 19.37985655913735, 19.98447548311175),
    @(16.99838664812632, 3.512343017039378, 23.33966680637614, 22.98623325781024, 3.582780986930106, 8.927716179194981, 15.39423010014013, 19.39419018621083, 19.70023522469667),
    @(16.65049104411149, 3.492365626306642, 23.03414365050934, 22.42546288780735, 3.5871775035615, 8.725720806751248, 15.20748146700314, 19.41868848295987, 19.41208297107439)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $rowNum = $i + 2
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($rowNum, $cols[$j]).Value = $rowValues[$j]
    }
}
